$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: new iteration V1.3.1 ---
$ws.Range("A8").Value = "V1.3.1"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "V1.3"
$ws.Range("D8").Value = 16
$ws.Range("E8").Value = "Adam"
$ws.Range("F8").Value = "Cross Entropy "
$ws.Range("G8").Value = 0.00000125
$ws.Range("G8").NumberFormat = "0.00E+00"
$ws.Range("K8").Value = 99.63
$ws.Range("L8").Value = 99.68
$ws.Range("N8").Value = "Took v1.3 classifier, retrained it decreasing the learning rate and result obtained is remarkable."

# --- Row 9: new iteration V1.3.2 ---
$ws.Range("A9").Value = "V1.3.2"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = "V1.3"
$ws.Range("D9").Value = 16
$ws.Range("E9").Value = "Adam"
$ws.Range("F9").Value = "Cross Entropy "
$ws.Range("G9").Value = 0.00000125
$ws.Range("G9").NumberFormat = "0.00E+00"
$ws.Range("K9").Value = 99.7
$ws.Range("L9").Value = 99.72

# --- Column N is now wider to fit the longer remark text ---
$ws.Columns.Item(14).ColumnWidth = 76.65

# --- Selection / view moved towards the newly added data ---
$ws.Range("N19").Select()
